# heterogeneity model with random slope
# Adds 8 new timeline rows (11-18) and re-applies "wrap text" formatting
# across the sheet (columns, header, body) to match the edited workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows -------------------------------------------------------
$ws.Range("A11").Value = "5th feb 2pm"
$ws.Range("B11").Value = "the issue of ones trick not working in jags posted on JAGS forum"

$ws.Range("A12").Value = "5 feb 5pm"
$ws.Range("B12").Value = "heterogeneity model can now find intercept easily. Constrained the sum of means to be zero"

$ws.Range("A13").Value = "5 feb 6:30pm"
$ws.Range("B13").Value = "step(abs(mu[1] + mu[2] + mu[3]) - 0.05) doesn" + [char]0x2019 + "t work as a good constraint. Really small values are sampled. Probably it is a better idea to freely sample params and then later transform them so that we get what we want."

$ws.Range("A14").Value = "6 feb 11pm"
$ws.Range("B14").Value = "in case of disproportionate proportions, even a simple GMM is difficult to estimate. It is tough to have good starting values for means but easier to do so with proportions. I think having highly disproportionate proportions at the beginning is a good idea."

$ws.Range("A15").Value = "6 feb 11:30pm"
$ws.Range("B15").Value = "Hmm probably try different proportions and see what results you get" + [char]0x2026 + "if less label switches then choose that one."

$ws.Range("A16").Value = "7 Feb noon"
$ws.Range("B16").Value = "choice of prior parameters for eta is tough especially when some of the components are nearly empty"

$ws.Range("A17").Value = "8 Feb noon"
$ws.Range("B17").Value = "Model with random slope running. Slow mixing and thus many estimates are quite poor."

$ws.Range("A18").Value = "8 Feb 1pm"
$ws.Range("B18").Value = "Model with random slope but effects were mean centred so that fixed effects could also be estimated. Working quite well it seems, some runs were bad otherwise all parameters are estimated well"

# --- Wrap text everywhere -------------------------------------------------
# The whole used range (header + all data rows) ends up with WrapText on.
$ws.Range("A1:B18").WrapText = $true

# --- Row heights for the new rows (auto-fit-ish heights from the source) --
$ws.Rows(12).RowHeight = 30
$ws.Rows(13).RowHeight = 45
$ws.Rows(14).RowHeight = 60
$ws.Rows(15).RowHeight = 30
$ws.Rows(16).RowHeight = 30
$ws.Rows(17).RowHeight = 30
$ws.Rows(18).RowHeight = 45

# A14 picked up an explicit (text) date format in the source workbook.
$ws.Range("A14").NumberFormat = "d-mmm"

# --- View state: scrolled down to the new rows, selection on B19 ---------
$ws.Range("B19").Select()
$ws.Application.ActiveWindow.ScrollRow = 13
